# Applies the "VTC Test Marking Scheme" content refresh:
#  - Prepends a short topical prefix to each question's text (B column) on
#    both the "Marking Scheme" and "Question Overview" sheets.
#  - Restyles the "### General Grading Guide (0-10 Scale)" markdown heading
#    to bold ("**General Grading Guide (0-10 Scale)**") inside the marking
#    scheme text (C column) on the "Marking Scheme" sheet.
#  - Refreshes the "Generated On" timestamp on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# Question-number -> short topical prefix, applied to rows 2..6 (Q1..Q5).
$prefixes = @{
    2 = "The Role of VTC: "
    3 = "Member Institutions: "
    4 = "Educational Philosophy: "
    5 = "Study Pathways: "
    6 = "Industry Partnership: "
}

$oldHeading = "### General Grading Guide (0-10 Scale)"
$newHeading = "**General Grading Guide (0-10 Scale)**"

# --- "Marking Scheme" sheet: question text prefixes + grading-guide heading ---
$wsScheme = $wb.Worksheets.Item("Marking Scheme")
foreach ($row in 2..6) {
    $bCell = $wsScheme.Cells.Item($row, 2)
    $bText = $bCell.Value()
    $bCell.Value = $prefixes[$row] + $bText

    $cCell = $wsScheme.Cells.Item($row, 3)
    $cText = $cCell.Value()
    $cCell.Value = $cText.Replace($oldHeading, $newHeading)
}

# --- "Question Overview" sheet: question text prefixes only ---
$wsOverview = $wb.Worksheets.Item("Question Overview")
foreach ($row in 2..6) {
    $bCell = $wsOverview.Cells.Item($row, 2)
    $bText = $bCell.Value()
    $bCell.Value = $prefixes[$row] + $bText
}

# --- "Summary" sheet: refresh the generated-on timestamp ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(7, 2).Value = "2026-01-05 01:23:05"
